# Commit: "Add ARM C and D"
# TC01 -> TC02 test-case rename: update the two generated workbook-name
# strings on the "startup" sheet, and move the active selection from E3
# (WebExcel filename cell) to D3 (TsvExcel filename cell) to match the
# saved cursor position in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the referenced TSV/Web data workbook file names for this test case.
$ws.Range("D2").Value = "TC02_Bento_Filter_Arm-B_TSVData.xlsx"
$ws.Range("E2").Value = "TC02_Bento_Filter_Arm-B_WebData.xlsx"

# Try to nudge the window position to match the author's saved view.
# (Best-effort; some hosts don't persist window geometry.)
try { $wb.Windows.Item(1).Left = 420 } catch {}
try { $wb.Windows.Item(1).Top = 5180 } catch {}

# Move the active cell/selection to D3, matching the saved workbook state.
$ws.Range("D3").Select()
